$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 9 ("「舌」لِسَانٌ ...") entirely; all rows below shift up by one.
$ws.Rows.Item(9).Delete()
